$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 2603
$ws.Range('I3').Value = 2737
$ws.Range('H4').Value = 1665
$ws.Range('I4').Value = 668
$ws.Range('I5').Value = 241
$ws.Range('I6').Value = 3132
$ws.Range('H7').Value = 25972
$ws.Range('I7').Value = 9381

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I2').Value = 27
$ws.Range('I6').Value = 34
$ws.Range('I7').Value = 100

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range('I3').Value = 2
$ws.Range('I6').Value = 4

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I2').Value = 30
$ws.Range('I6').Value = 47
$ws.Range('I7').Value = 114

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I2').Value = 98
$ws.Range('I3').Value = 94
$ws.Range('I7').Value = 308

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I3').Value = 60
$ws.Range('I6').Value = 49
$ws.Range('I7').Value = 171

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I2').Value = 81
$ws.Range('I3').Value = 125
$ws.Range('I6').Value = 126
$ws.Range('I7').Value = 358

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('I2').Value = 29
$ws.Range('I7').Value = 73

$ws = $wb.Worksheets.Item('New City')
$ws.Range('I2').Value = 64
$ws.Range('I7').Value = 210

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I5').Value = 30
$ws.Range('I6').Value = 65
$ws.Range('I7').Value = 310
$ws.Range('I8').Value = 595
$ws.Range('I11').Value = 157
$ws.Range('I19').Value = 261
$ws.Range('I22').Value = 26
$ws.Range('I25').Value = 43
$ws.Range('I27').Value = 89
$ws.Range('I29').Value = 626
$ws.Range('I33').Value = 443
$ws.Range('I37').Value = 308
$ws.Range('I38').Value = 4
$ws.Range('I41').Value = 45
$ws.Range('I42').Value = 325
$ws.Range('I47').Value = 66
$ws.Range('I48').Value = 102
$ws.Range('I51').Value = 81
$ws.Range('I52').Value = 191
$ws.Range('I53').Value = 104
$ws.Range('I54').Value = 209
$ws.Range('H55').Value = 310
$ws.Range('I55').Value = 100
$ws.Range('I59').Value = 18
$ws.Range('I63').Value = 34
$ws.Range('I65').Value = 210
$ws.Range('I66').Value = 22
$ws.Range('I67').Value = 358
$ws.Range('I69').Value = 24
$ws.Range('I71').Value = 21
$ws.Range('I76').Value = 146
$ws.Range('I78').Value = 128
$ws.Range('I79').Value = 242
$ws.Range('I83').Value = 183
$ws.Range('I84').Value = 73
$ws.Range('I85').Value = 435
$ws.Range('I86').Value = 54
$ws.Range('I88').Value = 79
$ws.Range('I89').Value = 100
$ws.Range('I91').Value = 110
$ws.Range('I92').Value = 28
$ws.Range('I96').Value = 114
$ws.Range('I97').Value = 72
$ws.Range('I99').Value = 171
$ws.Range('H101').Value = 25972
$ws.Range('I101').Value = 9381

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I2').Value = 64
$ws.Range('I6').Value = 34
$ws.Range('I7').Value = 183

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I3').Value = 161
$ws.Range('I5').Value = 11
$ws.Range('I6').Value = 144
$ws.Range('I7').Value = 443

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I3').Value = 42
$ws.Range('I7').Value = 209

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 195
$ws.Range('I3').Value = 219
$ws.Range('I5').Value = 21
$ws.Range('I6').Value = 170
$ws.Range('I7').Value = 626

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I3').Value = 70
$ws.Range('I7').Value = 261

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I2').Value = 13
$ws.Range('I3').Value = 20
$ws.Range('I6').Value = 58
$ws.Range('I7').Value = 102

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I2').Value = 31
$ws.Range('I3').Value = 35
$ws.Range('I6').Value = 60
$ws.Range('I7').Value = 146

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I2').Value = 108
$ws.Range('I3').Value = 178
$ws.Range('I4').Value = 21
$ws.Range('I7').Value = 435

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('I6').Value = 11
$ws.Range('I7').Value = 65

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('I3').Value = 14
$ws.Range('I7').Value = 45

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 83
$ws.Range('I3').Value = 111
$ws.Range('I5').Value = 13
$ws.Range('I6').Value = 88
$ws.Range('I7').Value = 325

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I2').Value = 24
$ws.Range('I3').Value = 32
$ws.Range('I6').Value = 51
$ws.Range('I7').Value = 128

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('I2').Value = 34
$ws.Range('H4').Value = 22
$ws.Range('I6').Value = 34
$ws.Range('H7').Value = 310
$ws.Range('I7').Value = 100

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('I3').Value = 4
$ws.Range('I7').Value = 24

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('I3').Value = 35
$ws.Range('I7').Value = 110

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I3').Value = 76
$ws.Range('I7').Value = 242

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I2').Value = 48
$ws.Range('I3').Value = 75
$ws.Range('I6').Value = 39
$ws.Range('I7').Value = 191

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('I3').Value = 13
$ws.Range('I7').Value = 43

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('I2').Value = 9
$ws.Range('I3').Value = 23
$ws.Range('I7').Value = 66

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('I2').Value = 6
$ws.Range('I6').Value = 10
$ws.Range('I7').Value = 22

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I3').Value = 33
$ws.Range('I4').Value = 14
$ws.Range('I7').Value = 157

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('I2').Value = 7
$ws.Range('I7').Value = 18

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I2').Value = 14
$ws.Range('I7').Value = 72

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('I6').Value = 12
$ws.Range('I7').Value = 28

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I3').Value = 22
$ws.Range('I6').Value = 30
$ws.Range('I7').Value = 79

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 191
$ws.Range('I3').Value = 162
$ws.Range('I6').Value = 186
$ws.Range('I7').Value = 595

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('I3').Value = 7
$ws.Range('I7').Value = 30

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('I2').Value = 21
$ws.Range('I7').Value = 89

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I3').Value = 4
$ws.Range('I7').Value = 54

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I2').Value = 14
$ws.Range('I6').Value = 40
$ws.Range('I7').Value = 81

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I3').Value = 30
$ws.Range('I7').Value = 104

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('I2').Value = 6
$ws.Range('I3').Value = 9
$ws.Range('I7').Value = 26

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('I2').Value = 6
$ws.Range('I7').Value = 21

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I2').Value = 106
$ws.Range('I3').Value = 93
$ws.Range('I6').Value = 79
$ws.Range('I7').Value = 310
